$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "hola prros"
$ws.Range("C3").Value = "juazjuaz"

$ws.Range("D4").Select()
